$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Withdrawal" (card id 11) to "Trinkets and Baubles" and buff its effect to gain 5 gold
$ws.Range("J13").Value = "Gain 5 gold."
$ws.Range("E13").Value = "Trinkets and Baubles"

# Fix spelling: "Eldrich Horror" -> "Eldritch Horror" (card id 14)
$ws.Range("E16").Value = "Eldritch Horror"

# Rename "Artillery" -> "Rain of Fire" (card id 31)
$ws.Range("E33").Value = "Rain of Fire"

# Rename "Reusable tools" -> "Desperate Times" (card id 15)
$ws.Range("E17").Value = "Desperate Times"

# Update the active cell selection to match the author's editing position
$ws.Range("E29").Select()

# Column E best-fit width recalculates because "Trinkets and Baubles" is now the longest name
$ws.Columns("E:E").ColumnWidth = 18.65
